{"js": "// Update the date line and all the two-digit-by-two-digit multiplication\n// answers in the practice sheet table to a new day's problem set.\n// Each entry is [oldText, newText]; values are unique in the document so\n// a plain exact-text search safely targets the right run.\nconst replacements = [\n  [\"2024-05-18 Saturday\", \"2024-05-19 Sunday\"],\n  [\"30\u00d719=570\", \"22\u00d782=1804\"],\n  [\"82\u00d749=4018\", \"27\u00d787=2349\"],\n  [\"49\u00d762=3038\", \"19\u00d742=798\"],\n  [\"66\u00d789=5874\", \"88\u00d732=2816\"],\n  [\"32\u00d735=1120\", \"91\u00d765=5915\"],\n  [\"52\u00d754=2808\", \"28\u00d768=1904\"],\n  [\"89\u00d730=2670\", \"47\u00d751=2397\"],\n  [\"39\u00d772=2808\", \"70\u00d727=1890\"],\n  [\"44\u00d725=1100\", \"82\u00d771=5822\"],\n  [\"88\u00d730=2640\", \"49\u00d786=4214\"],\n  [\"76\u00d797=7372\", \"79\u00d782=6478\"],\n  [\"26\u00d770=1820\", \"70\u00d787=6090\"],\n  [\"72\u00d714=1008\", \"74\u00d782=6068\"],\n  [\"50\u00d769=3450\", \"45\u00d758=2610\"],\n  [\"13\u00d747=611\", \"74\u00d788=6512\"],\n  [\"54\u00d788=4752\", \"55\u00d748=2640\"],\n  [\"11\u00d783=913\", \"59\u00d778=4602\"],\n  [\"82\u00d712=984\", \"58\u00d756=3248\"],\n  [\"55\u00d721=1155\", \"58\u00d758=3364\"],\n  [\"99\u00d723=2277\", \"77\u00d732=2464\"],\n  [\"42\u00d734=1428\", \"72\u00d736=2592\"],\n  [\"91\u00d726=2366\", \"33\u00d789=2937\"],\n  [\"38\u00d748=1824\", \"24\u00d794=2256\"],\n  [\"67\u00d784=5628\", \"53\u00d797=5141\"],\n  [\"32\u00d715=480\", \"41\u00d769=2829\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all the two-digit-by-two-digit multiplication\n# answers in the practice sheet table to a new day's problem set.\n# Each pair is unique in the document, so a plain whole-document\n# Find/Replace-all for each exact string safely targets the right run.\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @(\"2024-05-18 Saturday\", \"2024-05-19 Sunday\"),\n    @(\"30\u00d719=570\", \"22\u00d782=1804\"),\n    @(\"82\u00d749=4018\", \"27\u00d787=2349\"),\n    @(\"49\u00d762=3038\", \"19\u00d742=798\"),\n    @(\"66\u00d789=5874\", \"88\u00d732=2816\"),\n    @(\"32\u00d735=1120\", \"91\u00d765=5915\"),\n    @(\"52\u00d754=2808\", \"28\u00d768=1904\"),\n    @(\"89\u00d730=2670\", \"47\u00d751=2397\"),\n    @(\"39\u00d772=2808\", \"70\u00d727=1890\"),\n    @(\"44\u00d725=1100\", \"82\u00d771=5822\"),\n    @(\"88\u00d730=2640\", \"49\u00d786=4214\"),\n    @(\"76\u00d797=7372\", \"79\u00d782=6478\"),\n    @(\"26\u00d770=1820\", \"70\u00d787=6090\"),\n    @(\"72\u00d714=1008\", \"74\u00d782=6068\"),\n    @(\"50\u00d769=3450\", \"45\u00d758=2610\"),\n    @(\"13\u00d747=611\", \"74\u00d788=6512\"),\n    @(\"54\u00d788=4752\", \"55\u00d748=2640\"),\n    @(\"11\u00d783=913\", \"59\u00d778=4602\"),\n    @(\"82\u00d712=984\", \"58\u00d756=3248\"),\n    @(\"55\u00d721=1155\", \"58\u00d758=3364\"),\n    @(\"99\u00d723=2277\", \"77\u00d732=2464\"),\n    @(\"42\u00d734=1428\", \"72\u00d736=2592\"),\n    @(\"91\u00d726=2366\", \"33\u00d789=2937\"),\n    @(\"38\u00d748=1824\", \"24\u00d794=2256\"),\n    @(\"67\u00d784=5628\", \"53\u00d797=5141\"),\n    @(\"32\u00d715=480\", \"41\u00d769=2829\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}\n"}
